$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'basketball under pants',
    'softball gear for girls',
    'running capri',
    'softball compression sleeve',
    'youth softball compression sleeve',
    'running tights mens',
    'spandex men',
    'hockey kneepads',
    'padded leg sleeve',
    'mens basketball gear',
    'snowboarding padded shorts',
    'padded shorts snowboarding',
    'knee sleeve wrestling',
    'sleeve knee pads',
    'womens compression leggings',
    'airsoft knee pads',
    'mens compression tights 3 4',
    'basketball clothes for men',
    'men running tights',
    'knee pads nike',
    'knee pads mizuno',
    'knee pads bike',
    'yoga capri pants',
    'knee pads mtb',
    'knee pads skating',
    'mens workout tights',
    'mens basketball pants',
    'asics knee pads',
    'mens workout tights pants',
    'downhill knee pads',
    'men gym pants',
    'athletic capris',
    'valken knee pads',
    'woodland knee pads',
    'training tights men',
    'short tights for men',
    'ua compression pants',
    'men workout tights',
    'knee pads for exercise',
    'mens leggins',
    'nike kneepads',
    'youth football girdle with knee pads',
    'compression tights with pads',
    'compression leggings with knee pads',
    'basketball padded compression pants',
    'basketball knee pad pants',
    'padded compression pants men basketball',
    'padded tights men basketball',
    'tights with pads basketball',
    'mens leggings with knee pads',
    'leggings with knee pads women',
    'mtb knee pads men',
    'basketball padded knee sleeve',
    'compression pants women',
    'compression knee sleeve men basketball',
    'basketball sweat pants for men',
    'knee sleeve for wrestling',
    'leg sleeves for basketball youth',
    'training pants men',
    'compression knee sleeve men basketball',
    'basketball sweat pants for men',
    'knee sleeve for wrestling',
    'leg sleeves for basketball youth',
    'goalkeeper knee pads',
    'basketball calf sleeve',
    'compression knee sleeves with pads',
    'compression sleeve knee pads',
    'youth knee sleeve',
    'knee pad for scooter',
    'knee basketball',
    'knee pads for basketball youth',
    'mens compression knee',
    'knee pad sleeve basketball',
    'mens 3/4 compression pants',
    'youth compression knee pad sleeve',
    'mens basketball knee sleeves',
    'knee sleeve wrestling youth',
    'knee sleeves basketball men',
    'soccer compression pants',
    'leggings tight',
    'basketball leg sleeve youth padded',
    'knee pad construction',
    'youth basketball knee sleeve',
    'working knee pads for men',
    'cycling pants for men',
    'boys youth leggings',
    'compression running capris',
    'knee sleeve baseball',
    'compression knee sleeves for basketball',
    'volleyball kneepads',
    'compression knee sleeve with pad',
    'men capri shorts',
    'running compression pants',
    'mens work pants knee pads',
    'best knee pads',
    'compression pants sleeves',
    'mens compression running tights',
    'knee pads working',
    'basketball aids',
    'baseball youth pants'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
